$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Job# values in column B (B2:B4), keeping them as text and
# matching the border/alignment style already used by column A (style index 1). ---

# B2: "32297400" stays the same text, but needs the bordered style.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "32297400"
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)   # xlPasteFormats

# B3: "32297401" stays the same text, but needs the bordered style.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "32297401"
$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)   # xlPasteFormats

# B4: value actually changes to "32339569", needs the bordered style.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "32339569"
$ws.Range("A4").Copy()
$ws.Range("B4").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Remove the trailing blank row (row 5) ---
$ws.Rows("5").Delete()

# --- Update the active selection shown when the file is reopened ---
$ws.Range("E12").Select()
